$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update conversion rates embedded in the A1 text ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.55 = 17909.22 pesos`n✅ 17909.22 pesos = 4.52 = 936.11 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 219.999
$ws2.Range("O10").Value = 3940.01
$ws2.Range("N12").Value = 3961.94
$ws2.Range("O12").Value = 207.089
